$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: only date changes
$ws.Range("D2").Value = 44452

# Row 3: only date changes
$ws.Range("D3").Value = 44446

# Row 4: date + quality/volume/prices
$ws.Range("D4").Value = 44461
$ws.Range("L4").Value = "Especial"
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 31000
$ws.Range("O4").Value = 32000
$ws.Range("P4").Value = 31500
$ws.Range("S4").Value = 3150

# Row 5: date + quality/volume/prices
$ws.Range("D5").Value = 44461
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 30
$ws.Range("N5").Value = 30000
$ws.Range("O5").Value = 30000
$ws.Range("P5").Value = 30000
$ws.Range("S5").Value = 3000

# Row 6: date + volume/prices
$ws.Range("D6").Value = 44448
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 21000
$ws.Range("O6").Value = 22000
$ws.Range("P6").Value = 21500
$ws.Range("S6").Value = 2150

# Row 7: date + quality/prices
$ws.Range("D7").Value = 44447
$ws.Range("L7").Value = "Primera"
$ws.Range("N7").Value = 21000
$ws.Range("O7").Value = 22000
$ws.Range("P7").Value = 21500
$ws.Range("S7").Value = 2150

# Row 8: date + prices
$ws.Range("D8").Value = 44487
$ws.Range("N8").Value = 23000
$ws.Range("O8").Value = 24000
$ws.Range("P8").Value = 23500
$ws.Range("S8").Value = 2350

# Row 9: date + quality/prices
$ws.Range("D9").Value = 44460
$ws.Range("L9").Value = "Especial"
$ws.Range("N9").Value = 31000
$ws.Range("O9").Value = 32000
$ws.Range("P9").Value = 31500
$ws.Range("S9").Value = 3150

# Row 10: date + volume/prices
$ws.Range("D10").Value = 44460
$ws.Range("M10").Value = 30
$ws.Range("N10").Value = 30000
$ws.Range("O10").Value = 30000
$ws.Range("P10").Value = 30000
$ws.Range("S10").Value = 3000
